$wb = $excel.ActiveWorkbook

# This script applies the value updates captured in the commit diff for
# "Golem_Profits" (FFXIV crafting-profit tracker) across all 8 job sheets.
# Each row holds market/leve pricing data in columns H:N; the scheduled
# runner refreshed these numbers. Cells that the diff removes entirely
# (no longer populated) are cleared rather than zeroed.

# --- Sheet: ALC ---
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H15").Value = 617.6
$ws.Range("I15").Value = 617.6
$ws.Range("K15").Value = 1852.8
$ws.Range("M15").Value = -1683.8
$ws.Range("H28").Value = 523.5
$ws.Range("I28").Value = 624.5714
$ws.Range("K28").Value = 624.5714
$ws.Range("M28").Value = -139.5714
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()  # was -1138
$ws.Range("H92").Value = 838.8823
$ws.Range("I92").Value = 628.4167
$ws.Range("J92").Value = 1344
$ws.Range("K92").Value = 628.4167
$ws.Range("L92").Value = 1344
$ws.Range("M92").Value = 619.5833
$ws.Range("N92").Value = -3840
$ws.Range("H98").Value = 32500
$ws.Range("I98").Value = 30000
$ws.Range("J98").Value = 40000
$ws.Range("K98").Value = 30000
$ws.Range("L98").Value = 40000
$ws.Range("M98").Value = -28502
$ws.Range("N98").Value = -42996
$ws.Range("H106").Value = 1900
$ws.Range("I106").Value = 1900
$ws.Range("K106").Value = 1900
$ws.Range("M106").Value = -1269
$ws.Range("H113").Value = 2269.8572
$ws.Range("I113").Value = 2269.8572
$ws.Range("K113").Value = 2269.8572
$ws.Range("M113").Value = 984.1428000000001
$ws.Range("H122").Value = 32500
$ws.Range("I122").Value = 30000
$ws.Range("J122").Value = 40000
$ws.Range("K122").Value = 90000
$ws.Range("L122").Value = 120000
$ws.Range("M122").Value = -87550
$ws.Range("N122").Value = -124900

# --- Sheet: ARM ---
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 731.5833
$ws.Range("I32").Value = 681.2727
$ws.Range("K32").Value = 681.2727
$ws.Range("M32").Value = -394.2727
$ws.Range("H41").Value = 2397.8
$ws.Range("I41").Value = 2397.8
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 2397.8
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -1983.8
$ws.Range("N41").ClearContents()  # was -34828
$ws.Range("H61").Value = 2186.6667
$ws.Range("I61").Value = 2000
$ws.Range("K61").Value = 2000
$ws.Range("M61").Value = -1788
$ws.Range("H88").Value = 3995.4
$ws.Range("J88").Value = 3995.4
$ws.Range("L88").Value = 3995.4
$ws.Range("N88").Value = -4807.4
$ws.Range("H91").Value = 3995.4
$ws.Range("J91").Value = 3995.4
$ws.Range("L91").Value = 3995.4
$ws.Range("N91").Value = -6803.4
$ws.Range("H136").Value = 2186.6667
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

# --- Sheet: BSM ---
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H22").Value = 711.25
$ws.Range("I22").Value = 711.25
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 711.25
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -538.25
$ws.Range("N22").ClearContents()  # was -1148
$ws.Range("H94").Value = 2143.125
$ws.Range("I94").Value = 1798.75
$ws.Range("J94").Value = 2487.5
$ws.Range("K94").Value = 1798.75
$ws.Range("L94").Value = 2487.5
$ws.Range("M94").Value = -1347.75
$ws.Range("N94").Value = -3389.5
$ws.Range("H107").Value = 1732.7931
$ws.Range("I107").Value = 1560.3182
$ws.Range("K107").Value = 1560.3182
$ws.Range("M107").Value = 359.6818000000001
$ws.Range("H124").Value = 14307142
$ws.Range("I124").Value = 25007500
$ws.Range("J124").Value = 39999
$ws.Range("K124").Value = 25007500
$ws.Range("L124").Value = 39999
$ws.Range("M124").Value = -25002590
$ws.Range("N124").Value = -49819

# --- Sheet: CRP ---
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H7").Value = 74
$ws.Range("I7").Value = 72
$ws.Range("J7").Value = 78
$ws.Range("K7").Value = 72
$ws.Range("L7").Value = 78
$ws.Range("M7").Value = 41
$ws.Range("N7").Value = -304
$ws.Range("H15").Value = 10721.75
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 10721.75
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 10721.75
$ws.Range("M15").ClearContents()  # was -353.5294
$ws.Range("N15").Value = -11061.75
$ws.Range("H16").Value = 1172.75
$ws.Range("I16").Value = 1172.75
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1172.75
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -885.75
$ws.Range("N16").ClearContents()  # was -1575
$ws.Range("H22").Value = 734.8570999999999
$ws.Range("I22").Value = 752.9231
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 752.9231
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -402.9231
$ws.Range("N22").Value = -1200
$ws.Range("H31").Value = 1204.3334
$ws.Range("I31").Value = 1204.3334
$ws.Range("K31").Value = 1204.3334
$ws.Range("M31").Value = -909.3334
$ws.Range("H34").Value = 1204.3334
$ws.Range("I34").Value = 1204.3334
$ws.Range("K34").Value = 1204.3334
$ws.Range("M34").Value = -1002.3334
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()  # was -1702.8334
$ws.Range("H86").Value = 9069.143
$ws.Range("I86").Value = 9332.333000000001
$ws.Range("K86").Value = 9332.333000000001
$ws.Range("M86").Value = -8209.333000000001
$ws.Range("H89").Value = 9069.143
$ws.Range("I89").Value = 9332.333000000001
$ws.Range("K89").Value = 46661.665
$ws.Range("M89").Value = -41045.665
$ws.Range("H94").Value = 1392.5
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 1785
$ws.Range("K94").Value = 1000
$ws.Range("L94").Value = 1785
$ws.Range("M94").Value = -549
$ws.Range("N94").Value = -2687
$ws.Range("H105").Value = 714
$ws.Range("I105").Value = 666.3333
$ws.Range("K105").Value = 666.3333
$ws.Range("M105").Value = 1080.6667
$ws.Range("H113").Value = 1172.75
$ws.Range("I113").Value = 1172.75
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1172.75
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 997.25
$ws.Range("N113").ClearContents()  # was -5341
$ws.Range("H122").Value = 1946.8
$ws.Range("I122").Value = 1183.5
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 3550.5
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -1100.5
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 1391.6
$ws.Range("I132").Value = 1391.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4174.799999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1644.799999999999
$ws.Range("N132").ClearContents()  # was -10610
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()  # was -3167.5002

# --- Sheet: CUL ---
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H138").Value = 3998
$ws.Range("I138").Value = 3998
$ws.Range("K138").Value = 11994
$ws.Range("M138").Value = -6854
$ws.Range("H139").Value = 275
$ws.Range("I139").Value = 300
$ws.Range("K139").Value = 900
$ws.Range("M139").Value = 4240
$ws.Range("H140").Value = 230
$ws.Range("I140").Value = 230
$ws.Range("K140").Value = 690
$ws.Range("M140").Value = 4490

# --- Sheet: GSM ---
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H2").Value = 94.90909000000001
$ws.Range("I2").Value = 88.40000000000001
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 88.40000000000001
$ws.Range("L2").Value = 160
$ws.Range("M2").Value = 24.59999999999999
$ws.Range("N2").Value = -386
$ws.Range("H68").Value = 52000
$ws.Range("J68").Value = 52000
$ws.Range("L68").Value = 52000
$ws.Range("N68").Value = -53622
$ws.Range("H71").Value = 52000
$ws.Range("J71").Value = 52000
$ws.Range("L71").Value = 156000
$ws.Range("N71").Value = -164112
$ws.Range("H132").Value = 770.6667
$ws.Range("I132").Value = 770.6667
$ws.Range("K132").Value = 2312.0001
$ws.Range("M132").Value = 217.9998999999998

# --- Sheet: LTW ---
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H46").Value = 340332.34
$ws.Range("I46").Value = 2000000
$ws.Range("J46").Value = 8398.799999999999
$ws.Range("K46").Value = 2000000
$ws.Range("L46").Value = 8398.799999999999
$ws.Range("M46").Value = -1999812
$ws.Range("N46").Value = -8774.799999999999
$ws.Range("H62").Value = 50249
$ws.Range("J62").Value = 50249
$ws.Range("L62").Value = 50249
$ws.Range("N62").Value = -51497
$ws.Range("H65").Value = 50249
$ws.Range("J65").Value = 50249
$ws.Range("L65").Value = 150747
$ws.Range("N65").Value = -156987

# --- Sheet: WVR ---
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H4").Value = 1375.75
$ws.Range("J4").Value = 667.6667
$ws.Range("L4").Value = 667.6667
$ws.Range("N4").Value = -893.6667
$ws.Range("H64").Value = 10526
$ws.Range("J64").Value = 10526
$ws.Range("L64").Value = 10526
$ws.Range("N64").Value = -11022
$ws.Range("H67").Value = 10526
$ws.Range("J67").Value = 10526
$ws.Range("L67").Value = 10526
$ws.Range("N67").Value = -12242
$ws.Range("H81").Value = 2124.3333
$ws.Range("I81").Value = 2049.2
$ws.Range("K81").Value = 4098.4
$ws.Range("M81").Value = -3037.4
$ws.Range("H84").Value = 2124.3333
$ws.Range("I84").Value = 2049.2
$ws.Range("K84").Value = 20492
$ws.Range("M84").Value = -15188
$ws.Range("H132").Value = 1803.2727
$ws.Range("I132").Value = 1733.7
$ws.Range("K132").Value = 5201.1
$ws.Range("M132").Value = -2671.1
$ws.Range("H136").Value = 5051.364
$ws.Range("I136").Value = 5498.5
$ws.Range("J136").Value = 580
$ws.Range("K136").Value = 16495.5
$ws.Range("L136").Value = 1740
$ws.Range("M136").Value = -13945.5
$ws.Range("N136").Value = -6840
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()  # was -130280
